$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price refresh diff.
# Every cell in the sheet is stored as text (inline/shared string) in the
# source workbook, including Price values that look numeric (e.g. '9.400',
# '0.7645'). A bare Value assignment would let Excel auto-convert those to
# real numbers (losing trailing zeros / exact formatting), so any new value
# that parses as a number is written to a cell pre-formatted as Text ('@')
# to force it to round-trip as the literal string.

# Row 2
$ws.Range("D2").Value = '29.897.55'
$ws.Range("E2").Value = '  +0.06%  '
# Row 3
$ws.Range("D3").Value = '1.888.71'
$ws.Range("E3").Value = '  -0.20%  '
# Row 4
$ws.Range("E4").Value = '  +0.06%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7645'
$ws.Range("E5").Value = '  -1.55%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.78'
$ws.Range("E6").Value = '  -0.78%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
# Row 8
$ws.Range("E8").Value = '  -0.67%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.64'
$ws.Range("E9").Value = '  +0.92%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07166'
$ws.Range("E10").Value = '  -2.95%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08512'
$ws.Range("E11").Value = '  +4.62%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7623'
$ws.Range("E12").Value = '  -0.65%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.360'
$ws.Range("E13").Value = '  -2.34%  '
# Row 14
$ws.Range("D14").Value = '1.849.89'
$ws.Range("E14").Value = '  -1.92%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.64'
$ws.Range("E15").Value = '  +1.24%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.135'
$ws.Range("E16").Value = '  -1.39%  '
# Row 17
$ws.Range("D17").Value = '29.706.64'
$ws.Range("E17").Value = '  -0.58%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.77'
$ws.Range("E18").Value = '  -1.59%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.83'
$ws.Range("E19").Value = '  -0.48%  '
# Row 20
$ws.Range("E20").Value = '  -0.75%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9991'
$ws.Range("E21").Value = '  -0.14%  '
# Row 22
$ws.Range("D22").Value = '2.104.55'
$ws.Range("E22").Value = '  -2.20%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.958'
$ws.Range("E23").Value = '  -2.51%  '
# Row 24
$ws.Range("E24").Value = '  +0.04%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1618'
$ws.Range("E25").Value = '  +2.24%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.400'
$ws.Range("E26").Value = '  -0.53%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.89'
$ws.Range("E27").Value = '  -0.46%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("E28").Value = '  -0.51%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.034'
$ws.Range("E29").Value = '  -0.63%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.468'
$ws.Range("E30").Value = '  +1.32%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.535'
$ws.Range("E31").Value = '  -0.87%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.491'
$ws.Range("E32").Value = '  -0.47%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.094'
$ws.Range("E33").Value = '  -0.23%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05445'
$ws.Range("E34").Value = '  -2.85%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.242'
$ws.Range("E35").Value = '  -0.83%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7416'
$ws.Range("E36").Value = '  -3.08%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").Value = '  -0.07%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +1.91%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("E39").Value = '  +0.68%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.780'
$ws.Range("E40").Value = '  -0.42%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4463'
# Row 42
$ws.Range("D42").Value = '1.099.78'
$ws.Range("E42").Value = '  -5.25%  '
# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.063'
$ws.Range("E43").Value = '  +1.33%  '
# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.85'
$ws.Range("E44").Value = '  -2.13%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8531'
$ws.Range("E45").Value = '  -0.20%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.71'
$ws.Range("E47").Value = '  +0.57%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.865'
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.644'
$ws.Range("E49").Value = '  +1.27%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.996'
$ws.Range("E50").Value = '  -5.70%  '
# Row 51
$ws.Range("D51").Value = '2.015.62'
$ws.Range("E51").Value = '  -1.35%  '
